$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text (not numeric/percentage) interpretation for the edited range,
# matching the original inlineStr cell types, then reset style to avoid
# leaving a quote-prefix style applied.
$editRange = $ws.Range("D2:G51")
$editRange.NumberFormat = "@"

$ws.Range("D2").Value = "325.05"
$ws.Range("E2").Value = "-1.34%"
$ws.Range("G2").Value = "10"
$ws.Range("D3").Value = "39.51"
$ws.Range("E3").Value = "-1.43%"
$ws.Range("G3").Value = "10"
$ws.Range("D4").Value = "5.674"
$ws.Range("E4").Value = "7.22%"
$ws.Range("G4").Value = "10"
$ws.Range("D5").Value = "0.08015"
$ws.Range("E5").Value = "-1.13%"
$ws.Range("G5").Value = "10"
$ws.Range("D6").Value = "2.047"
$ws.Range("E6").Value = "6.25%"
$ws.Range("G6").Value = "10"
$ws.Range("D7").Value = "4.487"
$ws.Range("E7").Value = "-0.55%"
$ws.Range("G7").Value = "10"
$ws.Range("D8").Value = "8.620"
$ws.Range("E8").Value = "-0.23%"
$ws.Range("G8").Value = "10"
$ws.Range("D9").Value = "2.934"
$ws.Range("E9").Value = "-1.19%"
$ws.Range("G9").Value = "10"
$ws.Range("D10").Value = "0.9240"
$ws.Range("E10").Value = "-1.41%"
$ws.Range("G10").Value = "10"
$ws.Range("D11").Value = "0.1238"
$ws.Range("E11").Value = "-8.45%"
$ws.Range("G11").Value = "10"
$ws.Range("D12").Value = "0.1974"
$ws.Range("E12").Value = "-0.31%"
$ws.Range("G12").Value = "10"
$ws.Range("D13").Value = "8.741"
$ws.Range("E13").Value = "21.79%"
$ws.Range("G13").Value = "10"
$ws.Range("D14").Value = "0.09256"
$ws.Range("E14").Value = "0.34%"
$ws.Range("G14").Value = "10"
$ws.Range("D15").Value = "0.03710"
$ws.Range("E15").Value = "3.35%"
$ws.Range("G15").Value = "10"
$ws.Range("E16").Value = "9.31%"
$ws.Range("G16").Value = "10"
$ws.Range("D17").Value = "0.001292"
$ws.Range("E17").Value = "-2.65%"
$ws.Range("G17").Value = "10"
$ws.Range("D18").Value = "0.006139"
$ws.Range("E18").Value = "-4.18%"
$ws.Range("G18").Value = "10"
$ws.Range("D19").Value = "3.348"
$ws.Range("E19").Value = "-0.52%"
$ws.Range("G19").Value = "10"
$ws.Range("D20").Value = "0.3480"
$ws.Range("E20").Value = "-0.98%"
$ws.Range("G20").Value = "10"
$ws.Range("D21").Value = "0.1353"
$ws.Range("E21").Value = "2.32%"
$ws.Range("G21").Value = "10"
$ws.Range("E22").Value = "-5.77%"
$ws.Range("G22").Value = "10"
$ws.Range("D23").Value = "0.04395"
$ws.Range("E23").Value = "-0.72%"
$ws.Range("G23").Value = "10"
$ws.Range("D24").Value = "0.001252"
$ws.Range("E24").Value = "2.68%"
$ws.Range("G24").Value = "10"
$ws.Range("D25").Value = "0.004624"
$ws.Range("E25").Value = "6.18%"
$ws.Range("G25").Value = "10"
$ws.Range("E26").Value = "-3.35%"
$ws.Range("G26").Value = "10"
$ws.Range("G27").Value = "10"
$ws.Range("G28").Value = "10"
$ws.Range("G29").Value = "10"
$ws.Range("G30").Value = "10"
$ws.Range("G31").Value = "10"
$ws.Range("G32").Value = "10"
$ws.Range("G33").Value = "10"
$ws.Range("G34").Value = "10"
$ws.Range("G35").Value = "10"
$ws.Range("G36").Value = "10"
$ws.Range("G37").Value = "10"
$ws.Range("G38").Value = "10"
$ws.Range("D39").Value = "0.02494"
$ws.Range("E39").Value = "0.52%"
$ws.Range("G39").Value = "10"
$ws.Range("E40").Value = "2.60%"
$ws.Range("G40").Value = "10"
$ws.Range("D41").Value = "0.007475"
$ws.Range("E41").Value = "-2.81%"
$ws.Range("G41").Value = "10"
$ws.Range("D42").Value = "0.009841"
$ws.Range("E42").Value = "8.27%"
$ws.Range("G42").Value = "10"
$ws.Range("D43").Value = "0.1405"
$ws.Range("E43").Value = "-1.65%"
$ws.Range("G43").Value = "10"
$ws.Range("E44").Value = "-2.06%"
$ws.Range("G44").Value = "10"
$ws.Range("D45").Value = "0.01109"
$ws.Range("E45").Value = "12.05%"
$ws.Range("G45").Value = "10"
$ws.Range("D46").Value = "0.00006736"
$ws.Range("E46").Value = "1.08%"
$ws.Range("G46").Value = "10"
$ws.Range("E47").Value = "0.04%"
$ws.Range("G47").Value = "10"
$ws.Range("E48").Value = "-11.14%"
$ws.Range("G48").Value = "10"
$ws.Range("D49").Value = "0.002280"
$ws.Range("E49").Value = "-5.00%"
$ws.Range("G49").Value = "10"
$ws.Range("E50").Value = "0.04%"
$ws.Range("G50").Value = "10"
$ws.Range("E51").Value = "0.04%"
$ws.Range("G51").Value = "10"

$editRange.Style = "Normal"
